# Applies the "system properties" addition described in the commit:
#   se agrega cambio al system properties
#
# 1) Inserts four new paragraphs (indent 708, "nfasis" character style)
#    right after the "...standalone\configuration\standalone.xml" line
#    and before the "...ambientes clusterizados..." bullet, describing
#    the <system-properties> block to add inside <server ...>.
# 2) Updates the cached PAGE field result in the footer from "1" to "4"
#    (the document now spans more pages after the new content).

$d = $word.ActiveDocument

# --- 1. Locate the insertion point -----------------------------------
# Anchor on the unique "standalone.xml" path text; insert right after it
# (i.e. just before the paragraph mark that ends that paragraph).
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$null = $anchor.Find.Execute("standalone\configuration\standalone.xml", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertPos = $anchor.End
$insertRange = $d.Range($insertPos, $insertPos)

# --- 2. Build the OOXML for the four new paragraphs -------------------
$fragment = @'
<w:p><w:pPr><w:ind w:left="708"/><w:rPr><w:rStyle w:val="nfasis"/><w:i w:val="0"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="nfasis"/><w:i w:val="0"/></w:rPr><w:t xml:space="preserve">Dentro del tag </w:t></w:r><w:r><w:rPr><w:rStyle w:val="nfasis"/><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="nfasis"/><w:b/><w:i w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;server xmlns="urn:jboss:domain:1.2"&gt;  </w:t></w:r><w:r><w:rPr><w:rStyle w:val="nfasis"/><w:i w:val="0"/></w:rPr><w:t>agregue el siguiente contenido:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="708"/><w:rPr><w:rStyle w:val="nfasis"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="nfasis"/></w:rPr><w:t>&lt;system-properties&gt;</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="708"/><w:rPr><w:rStyle w:val="nfasis"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="nfasis"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rStyle w:val="nfasis"/></w:rPr><w:tab/><w:t>&lt;property name="org.apache.tomcat.util.http.Parameters.MAX_COUNT" value="5000"/&gt;</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="708"/><w:rPr><w:rStyle w:val="nfasis"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="nfasis"/></w:rPr><w:t>&lt;/system-properties&gt;</w:t></w:r></w:p>
'@

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $fragment +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($packageXml)

# --- 3. Refresh the footer's cached page-number field ------------------
# The document grew by a page, so the PAGE field's cached text in the
# footer must move from "1" to "4" to stay consistent.
$section = $d.Sections.Item(1)
$footer = $section.Footers.Item(1)
$pageField = $footer.Range.Fields.Item(1)
$pageField.Result.Text = "4"

Write-Output "done"
